# DOM and Banner author ids added
# - Update cited_by_count (M3) for row 3 from 2 to 4
# - Swap rows 4 and 5 (the "Tala B. Shahin..." row and the
#   "David F. Butler, Jonathan Skibo..." row trade places)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update M3: cited_by_count 2 -> 4
$ws.Range("M3").Value2 = 4

# Swap the contents of row 4 and row 5 across columns A..Q.
# Columns other than E (a date-like text column) can be copied directly;
# column E needs to be forced to stay as text so Excel does not silently
# convert the "yyyy-mm-dd" string into a date serial number.
$textCols = @(5)  # column E is index 5 (A=1 ... Q=17)
$lastCol = 17

for ($c = 1; $c -le $lastCol; $c++) {
    $cell4 = $ws.Cells.Item(4, $c)
    $cell5 = $ws.Cells.Item(5, $c)

    $val4 = $cell4.Value2
    $val5 = $cell5.Value2

    if ($textCols -contains $c) {
        $cell4.NumberFormat = "@"
        $cell5.NumberFormat = "@"
    }

    $cell4.Value2 = $val5
    $cell5.Value2 = $val4

    if ($textCols -contains $c) {
        $cell4.Style = "Normal"
        $cell5.Style = "Normal"
    }
}
